$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: "...1" -> "...3"
$ws.Range("A2").Value = "judul 3"
$ws.Range("B2").Value = "deskripsi singkat 3"
$ws.Range("C2").Value = "link 3"

# Row 3: "...2" -> "...4"
$ws.Range("A3").Value = "Judul 4"
$ws.Range("B3").Value = "deskripsi singkat 4"
$ws.Range("C3").Value = "link 4"

# Update the active selection to match the saved file (was B4, now C4)
$ws.Range("C4").Select()
